# Weekly refresh of the "Achicoria" market-price rows.
# The underlying records (43 rows) get reshuffled onto new dates/row
# positions, and one additional weekly observation is appended as row 45.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant for every data row and are
# left untouched; only D (Fecha), J (Volumen), K/L/M (Precio min/max/
# promedio), O (Origen) and P (Precio $/Kg) vary per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D=44438; J=34; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=3; D=44698; J=34; K=6000; L=7000; M=6500; O="Provincia de Quillota"; P=406},
    @{Row=4; D=44474; J=52; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=5; D=44994; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=6; D=44403; J=43; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=7; D=44999; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=8; D=44953; J=90; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=9; D=44407; J=45; K=5500; L=6000; M=5744; O="Provincia de Quillota"; P=359},
    @{Row=10; D=44313; J=34; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=11; D=44966; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=12; D=44575; J=61; K=8000; L=8000; M=8000; O="Provincia de Quillota"; P=500},
    @{Row=13; D=44355; J=25; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=14; D=44341; J=51; K=5500; L=6000; M=5755; O="Provincia de Quillota"; P=360},
    @{Row=15; D=44467; J=52; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=16; D=44589; J=52; K=8000; L=8000; M=8000; O="Provincia de Quillota"; P=500},
    @{Row=17; D=44442; J=25; K=6000; L=7000; M=6480; O="Provincia de Quillota"; P=405},
    @{Row=18; D=44973; J=90; K=7000; L=8000; M=7500; O="Provincia de Quillota"; P=469},
    @{Row=19; D=44582; J=52; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=20; D=44573; J=34; K=8000; L=8000; M=8000; O="Provincia de Quillota"; P=500},
    @{Row=21; D=44952; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=22; D=44358; J=52; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=23; D=44308; J=70; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=24; D=44971; J=160; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=25; D=44932; J=70; K=6000; L=7000; M=6500; O="Provincia de Quillota"; P=406},
    @{Row=26; D=44715; J=70; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=27; D=44946; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=28; D=44967; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=29; D=44477; J=25; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=30; D=44330; J=120; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=31; D=44376; J=43; K=4500; L=5000; M=4756; O="Provincia de Quillota"; P=297},
    @{Row=32; D=44691; J=61; K=6000; L=7000; M=6508; O="Provincia de Quillota"; P=407},
    @{Row=33; D=44943; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=34; D=44350; J=25; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=35; D=44455; J=52; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344},
    @{Row=36; D=44782; J=70; K=6000; L=6000; M=6000; O="Región Metropolitana"; P=375},
    @{Row=37; D=44363; J=160; K=5500; L=6000; M=5750; O="Provincia de Quillota"; P=359},
    @{Row=38; D=44371; J=34; K=5500; L=6000; M=5750; O="Provincia de Quillota"; P=359},
    @{Row=39; D=44938; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=40; D=44910; J=70; K=6000; L=7000; M=6500; O="Provincia de Quillota"; P=406},
    @{Row=41; D=44328; J=160; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=42; D=44957; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=43; D=44939; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438},
    @{Row=44; D=44306; J=50; K=6000; L=6000; M=6000; O="Provincia de Quillota"; P=375},
    @{Row=45; D=44960; J=70; K=7000; L=7000; M=7000; O="Provincia de Quillota"; P=438}
)

foreach ($r in $rows) {
    $n = $r.Row
    if ($n -eq 45) {
        $ws.Range("A45").Value = 9
        $ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
        $ws.Range("C45").Value = "Metropolitana"
        $ws.Range("E45").Value = 13
        $ws.Range("F45").Value = 100112010
        $ws.Range("G45").Value = "Achicoria"
        $ws.Range("H45").Value = "Sin especificar"
        $ws.Range("I45").Value = "Primera"
        $ws.Range("N45").Value = "`$/caja 16 unidades"
        $ws.Range("Q45").Value = 16
        $ws.Range("R45").Value = "Hortaliza"
    }

    $ws.Range("D$n").Value = $r.D
    $ws.Range("D$n").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
}

Write-Output "done"
